# Rewrite the forecast data table (rows 2-19) with the corrected
# evaluation / simulated rt_data values for every component column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is brand new (the table grew from 18 to 19 data rows); give its
# date cell (column A) the same date number-format style as the other
# date cells by copying the format down from A18 before filling values.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 5.896808312953783
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 10.70171490310616

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 7.441962824572235
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 10.14224506046018

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = 6.277541464866987
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 9.926356894615008

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 6.535114773304773
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 6.670100100023779

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 5.12051970717502
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 6.878526556821973

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 3.65682115264816
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 4.783388268072009

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 2.943878639034381
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 5.440628060653574

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 1.172679597477866
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 0.7846309372076199

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 2.961845079861303
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 2.521390070591267

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 2.508469427909898
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 2.540874511056623

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 3.523703831572056
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 3.025650759929999

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 1.178605266817589
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 0.6761281928490348

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 3.047037961814492
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 3.638425237116749

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -0.2228847697281378
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = 2.646620531497335

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = -1.165854108406617
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 4.306257314632855

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 2.501311189006916
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = 3.687525625406263

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = 0.6753076481029074
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.5720316833079497

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = 2.039329803030121
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 2.838786065810939

